# Sierra Leone master data
# Adapt the status_type master-data sheet from the Madagascar (French) template
# to the Sierra Leone (English) template:
#   - lang_code: fra -> eng
#   - descr / name text: "Statut d'activation" -> "Activation Status"
#   - is_active: store the literal word TRUE as text (not a boolean)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# language code
$ws.Range("A2").Value = "eng"

# name / descr translated to English
$ws.Range("C2").Value = "Activation Status"
$ws.Range("D2").Value = "Activation Status"

# is_active stored as literal text "TRUE" (force text type with a leading
# apostrophe so Excel does not coerce it back into a boolean)
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "'TRUE"

# widen the name/descr columns a bit and give column E a narrow default width
$ws.Columns.Item(3).ColumnWidth = 22.6328125
$ws.Columns.Item(4).ColumnWidth = 19.90625

# a couple of left-aligned, word-wrapped blank cells were formatted below the
# table (row 4, columns C:D)
$ws.Range("C4:D4").HorizontalAlignment = -4131
$ws.Range("C4:D4").WrapText = $true

# leave the cursor where the author left it
$ws.Range("E10").Select() | Out-Null
